$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.816.39"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "3.464.26"
$ws.Range("E3").Value = "  +4.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Formula = "=""240.12"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").Formula = "=""645.06"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +6.54%  "
$ws.Range("E8").Value = "  -3.20%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("D11").Value = "3.461.36"
$ws.Range("E11").Value = "  +4.18%  "
$ws.Range("E12").Value = "  -3.09%  "
$ws.Range("D13").Formula = "=""41.99"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +4.72%  "
$ws.Range("D14").Formula = "=""6.14"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "94.677.15"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "4.113.09"
$ws.Range("E16").Value = "  +4.36%  "
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Formula = "=""8.53"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "3.455.94"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("D20").Formula = "=""17.92"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +5.82%  "
$ws.Range("D21").Formula = "=""11.40"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +9.33%  "
$ws.Range("D22").Formula = "=""0.514"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Formula = "=""502.38"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -5.28%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Formula = "=""6.64"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Formula = "=""91.98"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("D28").Formula = "=""12.17"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "3.648.89"
$ws.Range("E29").Value = "  +4.29%  "
$ws.Range("D30").Formula = "=""11.74"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +12.68%  "
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").Formula = "=""0.184"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("D35").Formula = "=""31.02"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +11.45%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").Formula = "=""7.78"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Formula = "=""532.91"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Formula = "=""0.928"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +12.44%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Formula = "=""0.151"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Formula = "=""24.09"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("E45").Value = "  +4.96%  "
$ws.Range("D46").Formula = "=""1.70"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +2.78%  "
$ws.Range("D47").Formula = "=""0.0417"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -2.18%  "
$ws.Range("D48").Formula = "=""3.52"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  +10.18%  "
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").Formula = "=""53.31"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.12%  "
$excel.CutCopyMode = 0
